{"js": "const body = context.document.body;\nconst paragraph = body.paragraphs.getFirst();\nparagraph.insertText(\"This is still burning.\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs.Item(1)\n$r = $d.Range($p.Range.Start, $p.Range.End - 1)\n$r.Text = \"This is still burning.\"\n"}
